# Updated cryptos list on Wed Oct 23 17:00:14 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "Price" (column D) cell's value while forcing text storage
# for values that would otherwise be auto-parsed by Excel as a number
# (losing trailing zeros / exact formatting, e.g. "1.00" -> 1).
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "65.864.96"
$ws.Range("E2").Value = "  -1.67%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.495.43"
$ws.Range("E3").Value = "  -4.37%  "

# Row 5 - BNB
Set-TextValue "D5" "575.46"
$ws.Range("E5").Value = "  -2.87%  "

# Row 6 - Solana
Set-TextValue "D6" "165.83"
$ws.Range("E6").Value = "  -0.21%  "

# Row 8 - XRP
Set-TextValue "D8" "0.519"
$ws.Range("E8").Value = "  -2.35%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.498.83"
$ws.Range("E9").Value = "  -4.21%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.136"
$ws.Range("E10").Value = "  -0.71%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.26%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.346"
$ws.Range("E12").Value = "  -4.01%  "

# Row 13 - Toncoin
Set-TextValue "D13" "5.08"
$ws.Range("E13").Value = "  -2.52%  "

# Row 14 - Avalanche
Set-TextValue "D14" "26.18"
$ws.Range("E14").Value = "  -4.75%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.959.92"
$ws.Range("E15").Value = "  -4.14%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0000174"
$ws.Range("E16").Value = "  -3.75%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "65.815.04"
$ws.Range("E17").Value = "  -1.69%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.429.05"
$ws.Range("E18").Value = "  -8.26%  "

# Row 19 - Chainlink
Set-TextValue "D19" "11.18"
$ws.Range("E19").Value = "  -6.56%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.59"
$ws.Range("E20").Value = "  -4.19%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "343.37"
$ws.Range("E21").Value = "  -3.30%  "

# Row 22 - Polkadot
Set-TextValue "D22" "4.17"
$ws.Range("E22").Value = "  -3.06%  "

# Row 23 - NEARProtocol
Set-TextValue "D23" "4.52"
$ws.Range("E23").Value = "  -2.41%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.01%  "

# Row 25 - SuiNetwork
$ws.Range("E25").Value = "  +0.67%  "

# Row 26 - Litecoin
Set-TextValue "D26" "68.70"
$ws.Range("E26").Value = "  -1.26%  "

# Row 27 - Aptos
Set-TextValue "D27" "9.88"
$ws.Range("E27").Value = "  -2.69%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("E28").Value = "  +0.11%  "

# Row 29 - WrappedeETH
$ws.Range("E29").Value = "  -3.83%  "

# Row 30 - PEPE
Set-TextValue "D30" "0.0₃0966"
$ws.Range("E30").Value = "  -2.83%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "8.04"
$ws.Range("E31").Value = "  +2.50%  "

# Row 32 - Bittensor
Set-TextValue "D32" "514.97"
$ws.Range("E32").Value = "  -4.75%  "

# Row 33 - Fetch.AI
Set-TextValue "D33" "1.29"
$ws.Range("E33").Value = "  -3.33%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  -4.91%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.130"
$ws.Range("E35").Value = "  -3.08%  "

# Row 36 - FirstDigitalUSD
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  +0.00%  "

# Row 37 - Monero
Set-TextValue "D37" "157.34"
$ws.Range("E37").Value = "  -0.28%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "1.43"
$ws.Range("E38").Value = "  -3.87%  "

# Row 39 - EthereumClassic
Set-TextValue "D39" "18.39"
$ws.Range("E39").Value = "  -2.76%  "

# Row 40 - WhiteBITCoin
$ws.Range("E40").Value = "  +0.71%  "

# Row 41 - PolygonEcosystemToken
Set-TextValue "D41" "0.351"
$ws.Range("E41").Value = "  -3.55%  "

# Row 42 - Stacks
Set-TextValue "D42" "1.74"
$ws.Range("E42").Value = "  -2.85%  "

# Row 43 - RenderToken
Set-TextValue "D43" "5.01"
$ws.Range("E43").Value = "  -3.07%  "

# Row 45 - dogwifhat
Set-TextValue "D45" "2.40"
$ws.Range("E45").Value = "  +0.50%  "

# Row 46 - Aave
Set-TextValue "D46" "145.93"
$ws.Range("E46").Value = "  -3.55%  "

# Row 47 - ARBITRUM
Set-TextValue "D47" "0.548"
$ws.Range("E47").Value = "  -4.54%  "

# Row 48 - now Filecoin (was BabyDogeCoin)
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D48" "3.67"
$ws.Range("E48").Value = "  -1.80%  "

# Row 49 - now BabyDogeCoin (was Filecoin)
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0273"
$ws.Range("E49").Value = "  -7.28%  "

# Row 50 - Optimism
Set-TextValue "D50" "1.69"
$ws.Range("E50").Value = "  +0.46%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -2.41%  "
